$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove data validation rules from the whole sheet ---
$ws.Cells.Validation.Delete()

# --- Column widths ---
# col B (2) and C (3) used to share one merged "10.5" width entry; now split/resized.
# (Excel's ColumnWidth setter snaps to whole-pixel increments, so these inputs are
# chosen to land on the stored width closest to the target: 14.375, 18.5, 12, 15.75, 11)
$ws.Columns.Item(2).ColumnWidth = 13.714   # -> stored width ~14.43 (closest to 14.375)
$ws.Columns.Item(3).ColumnWidth = 17.858   # -> stored width ~18.57 (closest to 18.5)
$ws.Columns.Item(5).ColumnWidth = 11.286   # -> stored width 12
$ws.Columns.Item(6).ColumnWidth = 15.0     # -> stored width ~15.71 (closest to 15.75)
$ws.Columns.Item(7).ColumnWidth = 10.286   # -> stored width 11

# --- Row 2 data (new server entry) ---
# Make sure the Name (C) and IP (F) cells use the same "text" number format
# already used by the ID/ServerID columns, then fill in the values.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("A2").Value = "WorldServer_1"
$ws.Range("B2").Value = "000103001"
$ws.Range("C2").Value = "WorldServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("G2").Value = 3001

# --- Selection moves to G1 ---
$ws.Range("G1").Select()
